$wb = $excel.ActiveWorkbook
$wsTypo = $wb.Worksheets.Item("Typography")
$wsTrans = $wb.Worksheets.Item("Translation")

# Typography sheet: I4 = Wildcard Ranges value for Default font
$wsTypo.Range("I4").Value = "0-9"

# Translation sheet: update rows for photoresistor (light sensor) support
$wsTrans.Range("B4").Value = "SingleUseId2"
$wsTrans.Range("C4").Value = "Default"
$wsTrans.Range("D4").Value = "Center"
$wsTrans.Range("E4").Value = "LTR"
$wsTrans.Range("F4").Value = "<"

$wsTrans.Range("B5").Value = "SingleUseId3"
$wsTrans.Range("C5").Value = "Default"
$wsTrans.Range("D5").Value = "Center"
$wsTrans.Range("E5").Value = "LTR"
$wsTrans.Range("F5").Value = ">"

$wsTrans.Range("B6").Value = "SingleUseId12"
$wsTrans.Range("C6").Value = "Default"
$wsTrans.Range("D6").Value = "Center"
$wsTrans.Range("E6").Value = "LTR"
$wsTrans.Range("F6").Value = "Temperature:"

$wsTrans.Range("B7").Value = "SingleUseId11"
$wsTrans.Range("C7").Value = "Default"
$wsTrans.Range("D7").Value = "Center"
$wsTrans.Range("E7").Value = "LTR"
$wsTrans.Range("F7").Value = ">"

$wsTrans.Range("B8").Value = "SingleUseId10"
$wsTrans.Range("C8").Value = "Default"
$wsTrans.Range("D8").Value = "Center"
$wsTrans.Range("E8").Value = "LTR"
$wsTrans.Range("F8").Value = "<"

$wsTrans.Range("B9").Value = "SingleUseId14"
$wsTrans.Range("C9").Value = "Default"
$wsTrans.Range("D9").Value = "Center"
$wsTrans.Range("E9").Value = "LTR"
$wsTrans.Range("F9").Value = "<"

$wsTrans.Range("B10").Value = "SingleUseId15"
$wsTrans.Range("C10").Value = "Default"
$wsTrans.Range("D10").Value = "Center"
$wsTrans.Range("E10").Value = "LTR"
$wsTrans.Range("F10").Value = ">"

$wsTrans.Range("B11").Value = "SingleUseId16"
$wsTrans.Range("C11").Value = "Default"
$wsTrans.Range("D11").Value = "Center"
$wsTrans.Range("E11").Value = "LTR"
$wsTrans.Range("F11").Value = "Light:"

$wsTrans.Range("B12").Value = "SingleUseId17"
$wsTrans.Range("C12").Value = "Default"
$wsTrans.Range("D12").Value = "Center"
$wsTrans.Range("E12").Value = "LTR"
$wsTrans.Range("F12").Value = "<value>%"

$wsTrans.Range("B13").Value = "SingleUseId18"
$wsTrans.Range("C13").Value = "Default"
$wsTrans.Range("D13").Value = "Center"
$wsTrans.Range("E13").Value = "LTR"
$wsTrans.Range("F13").Value = "<"

$wsTrans.Range("B14").Value = "SingleUseId19"
$wsTrans.Range("C14").Value = "Default"
$wsTrans.Range("D14").Value = "Center"
$wsTrans.Range("E14").Value = "LTR"
$wsTrans.Range("F14").Value = ">"

$wsTrans.Range("B15").Value = "SingleUseId20"
$wsTrans.Range("C15").Value = "Default"
$wsTrans.Range("D15").Value = "Left"
$wsTrans.Range("E15").Value = "LTR"
$wsTrans.Range("F15").Value = "Red:"

$wsTrans.Range("B16").Value = "SingleUseId21"
$wsTrans.Range("C16").Value = "Default"
$wsTrans.Range("D16").Value = "Left"
$wsTrans.Range("E16").Value = "LTR"
$wsTrans.Range("F16").Value = "Green:"

$wsTrans.Range("B17").Value = "SingleUseId22"
$wsTrans.Range("C17").Value = "Default"
$wsTrans.Range("D17").Value = "Left"
$wsTrans.Range("E17").Value = "LTR"
$wsTrans.Range("F17").Value = "Blue:"

$wsTrans.Range("B18").Value = "SingleUseId23"
$wsTrans.Range("C18").Value = "Default"
$wsTrans.Range("D18").Value = "Left"
$wsTrans.Range("E18").Value = "LTR"
$wsTrans.Range("F18").Value = "MQTT Messages:"

$wsTrans.Range("B19").Value = "SingleUseId25"
$wsTrans.Range("C19").Value = "Default"
$wsTrans.Range("D19").Value = "Left"
$wsTrans.Range("E19").Value = "LTR"
$wsTrans.Range("F19").Value = "<value>%"

$wsTrans.Range("B20").Value = "SingleUseId26"
$wsTrans.Range("C20").Value = "Default"
$wsTrans.Range("D20").Value = "Left"
$wsTrans.Range("E20").Value = "LTR"
$wsTrans.Range("F20").Value = "0"

$wsTrans.Range("B21").Value = "SingleUseId27"
$wsTrans.Range("C21").Value = "Default"
$wsTrans.Range("D21").Value = "Left"
$wsTrans.Range("E21").Value = "LTR"
$wsTrans.Range("F21").Value = "0"
